$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Population" label from J1 while keeping its formatting
$ws.Range("J1").ClearContents()

# Update selection to J1 as in the final saved state
$ws.Range("J1").Select()
